$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.988.37'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.236.77'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.39'
$ws.Range("E5").Value = '  -4.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.56'
$ws.Range("E6").Value = '  -6.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E9").Value = '  -4.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.75'
$ws.Range("E10").Value = '  -5.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0807'
$ws.Range("E11").Value = '  -2.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.20'
$ws.Range("E12").Value = '  -4.42%  '
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.577.90'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.236.66'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.825'
$ws.Range("E16").Value = '  -3.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.60'
$ws.Range("E17").Value = '  -4.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.855.11'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.10'
$ws.Range("E20").Value = '  -9.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.15'
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.93'
$ws.Range("E24").Value = '  -5.06%  '
$ws.Range("E25").Value = '  -5.75%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.99'
$ws.Range("E27").Value = '  -5.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.43'
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.86'
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.06'
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0800'
$ws.Range("E33").Value = '  -5.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.25'
$ws.Range("E34").Value = '  +5.03%  '
$ws.Range("E35").Value = '  -3.77%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.107'
$ws.Range("E37").Value = '  -5.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  -7.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.10'
$ws.Range("E39").Value = '  -7.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.84'
$ws.Range("E40").Value = '  -7.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("E41").Value = '  -9.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0300'
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.727.68'
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '85.37'
$ws.Range("E45").Value = '  +5.36%  '
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.89'
$ws.Range("E47").Value = '  -3.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.92'
$ws.Range("E48").Value = '  -4.82%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.09'
$ws.Range("E49").Value = '  -2.35%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.95'
$ws.Range("E50").Value = '  -7.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.10'
$ws.Range("E51").Value = '  -5.59%  '
